# Updating PM and Risks
# Risk likelihood values (column G) were revised downward for a handful of
# risk rows; the dependent "score" formula in column H (=F*G) recalculates
# automatically. Also moves the active-cell selection from I10 to I9.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Risks")

$ws.Range("G5").Value = 0.2
$ws.Range("G10").Value = 0.4
$ws.Range("G11").Value = 0.5
$ws.Range("G12").Value = 0.7
$ws.Range("G13").Value = 1

$ws.Range("I9").Select() | Out-Null
